# Updated loading_percent results for the "case with 380 kV done" run.
# New per-row values for columns B, C, D, F, G, J, K, M (row r corresponds to A{r} = r-2).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2 = @(14.47291449158059, 6.536581162841211, 4.871905231262534, 62.06092919674609, 3.774943278481869, 11.37081718558679, 14.14630398073385, 17.47532753946232)
    3 = @(14.42025179919829, 6.543801668682601, 4.906647019855798, 61.09484299526316, 3.778522432661807, 11.37113637327892, 14.14895675237627, 17.53448646930017)
    4 = @(14.39345028580644, 6.551971826141193, 4.930574101440207, 60.49682920091463, 3.780831702119414, 11.37264039938179, 14.15564004274797, 17.57564516552668)
    5 = @(14.3839293669531, 6.556235736842426, 4.940970284289214, 60.25212234567796, 3.781800934941888, 11.37358243093718, 14.15963293741469, 17.59363036454534)
    6 = @(14.38243327044548, 6.55700002664525, 4.942735357093821, 60.21143387813104, 3.781963580964586, 11.37375873775556, 14.16037255641582, 17.59668994364)
    7 = @(14.39331620032762, 6.552025554592651, 4.930711702827713, 60.49353282281835, 3.780844659260642, 11.37265177109543, 14.15568875550937, 17.57588281430838)
    8 = @(14.45361423406712, 6.538292486810468, 4.883341344816339, 61.72894722679776, 3.776154263598471, 11.37065590265125, 14.14616895276206, 17.49472026587704)
    9 = @(14.61524548673175, 6.541185153172907, 4.811349593128554, 64.10446112810982, 3.767837296718805, 11.37711008457564, 14.16762996161678, 17.37405388778605)
    10 = @(14.75960404285603, 6.561646049338546, 4.771613096195666, 65.80958602653023, 3.762256764934819, 11.38815622003698, 14.20782613780845, 17.30903485463241)
    11 = @(14.83061572284809, 6.574942512525666, 4.756478637759268, 66.57425848717847, 3.759831596688685, 11.39454525226673, 14.23138632445057, 17.2846197133164)
    12 = @(14.85825284947748, 6.580549678493816, 4.75117748462003, 66.86205392866384, 3.758929445314412, 11.39716014293696, 14.24106204397911, 17.27611887869597)
    13 = @(14.85226785933555, 6.579316668729711, 4.75229995355092, 66.80015363612588, 3.759123020616129, 11.39658829571906, 14.23894475460736, 17.27791653298031)
    14 = @(14.8328746139524, 6.575392365120199, 4.756033853322514, 66.59797210212297, 3.759757051941173, 11.39475646677278, 14.23216727051597, 17.28390540833398)
    15 = @(14.82109224220961, 6.573063051972626, 4.7583771665993, 66.47389419775001, 3.760147521954524, 11.39365985804841, 14.22811390421132, 17.28767080694354)
    16 = @(14.75506894851956, 6.5608572594197, 4.772661961100902, 65.75937533766238, 3.762417529103338, 11.38776607703885, 14.2063922053854, 17.31073452445612)
    17 = @(14.71591919692699, 6.554390354742734, 4.782183676432268, 65.3180882430062, 3.763839085179383, 11.38449937374224, 14.19441478141501, 17.32620729293517)
    18 = @(14.69390554152715, 6.551046368761051, 4.787936817349775, 65.06324722172897, 3.764667410799189, 11.38274889160878, 14.18802251010352, 17.33559260210113)
    19 = @(14.68653933331126, 6.549978687448715, 4.789932007971476, 64.97679235206506, 3.764949705865469, 11.38217828639574, 14.185943646267, 17.33885366538802)
    20 = @(14.72003471902847, 6.555039900714929, 4.781141406097806, 65.36517141134915, 3.763686653110364, 11.38483383136141, 14.19563840829929, 17.32450989632943)
    21 = @(14.83855080000483, 6.576529520244619, 4.754925394099653, 66.65740720050761, 3.75957038255473, 11.39528921907807, 14.23413756038672, 17.28212610375465)
    22 = @(14.92034817054367, 6.593907177252103, 4.740300768844425, 67.49156289446037, 3.756974580909084, 11.40326159507861, 14.26369092089191, 17.25876712034314)
    23 = @(14.87630181367616, 6.584328252335945, 4.747874403552987, 67.04736950517184, 3.758351404632819, 11.39890259130473, 14.2475176533095, 17.2708363259976)
    24 = @(14.71817255091198, 6.554745076139129, 4.781611748051464, 65.3438886448471, 3.763755533207252, 11.3846822256104, 14.1950836683863, 17.32527576379185)
    25 = @(14.5669538675286, 6.537189028919395, 4.828544475297902, 63.46806799704111, 3.769993683729399, 11.37425544107152, 14.15752731412132, 17.40255817757207)
}

$cols = @("B", "C", "D", "F", "G", "J", "K", "M")

foreach ($row in $newValues.Keys) {
    $values = $newValues[$row]
    for ($i = 0; $i -lt $cols.Count; $i++) {
        $ws.Range("$($cols[$i])$row").Value = $values[$i]
    }
}
